$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = 294

$ws.Range("D81").Value = 177
$ws.Range("E81").Value = 818

$ws.Range("D92").Value = 94
$ws.Range("E92").Value = 553

$ws.Range("B98").Value = 508
$ws.Range("C98").Value = 6
$ws.Range("D98").Value = 294
$ws.Range("E98").Value = 205
$ws.Range("F98").Value = 12

$ws.Range("C99").Value = 28
$ws.Range("D99").Value = 31
$ws.Range("E99").Value = 431
$ws.Range("F99").Value = 3
$ws.Range("H99").Value = 31

$ws.Range("B100").Value = 493
$ws.Range("D100").Value = 159
$ws.Range("E100").Value = 317
$ws.Range("F100").Value = 2
$ws.Range("H100").Value = 17

$ws.Range("B101").Value = 489
$ws.Range("D101").Value = 114
$ws.Range("E101").Value = 370
$ws.Range("F101").Value = 5
$ws.Range("H101").Value = 5

$ws.Range("B102").Value = 477
$ws.Range("D102").Value = 59
$ws.Range("E102").Value = 415
$ws.Range("F102").Value = 0
$ws.Range("H102").Value = 3

$ws.Range("D120").Value = 113
$ws.Range("E120").Value = 105

$ws.Range("B126").Value = 163
$ws.Range("C126").Value = 20
$ws.Range("D126").Value = 25
$ws.Range("E126").Value = 133
$ws.Range("F126").Value = 0
$ws.Range("H126").Value = 5

$ws.Range("B127").Value = 158
$ws.Range("D127").Value = 73
$ws.Range("E127").Value = 77
$ws.Range("F127").Value = 17
$ws.Range("H127").Value = 8

$ws.Range("B128").Value = 147
$ws.Range("D128").Value = 11
$ws.Range("E128").Value = 131
$ws.Range("F128").Value = 4
$ws.Range("H128").Value = 5

$ws.Range("B129").Value = 145
$ws.Range("D129").Value = 67
$ws.Range("E129").Value = 70
$ws.Range("F129").Value = 13
$ws.Range("H129").Value = 8

$ws.Range("D130").Value = 11
$ws.Range("E130").Value = 126
$ws.Range("H130").Value = 6

$ws.Range("D155").Value = 2
$ws.Range("E155").Value = 53

$ws.Range("D195").Value = 0
$ws.Range("F195").Value = 0
$ws.Range("H195").Value = 1

$ws.Range("F196").Value = 1

$ws.Range("D197").Value = 1
$ws.Range("H197").Value = 0
